$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 162691
$ws.Range("C4").Value = 153689
$ws.Range("C5").Value = 9002
$ws.Range("C8").Value = 64.68000000000001
